# Update gh-pages output data (想去人数 / 最低票价 columns) across sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 218
$ws.Range("F3").Value = 54768
$ws.Range("F4").Value = 1344
$ws.Range("F5").Value = 385
$ws.Range("F6").Value = 327
$ws.Range("F7").Value = 881
$ws.Range("G7").Value = 60
$ws.Range("F8").Value = 760
$ws.Range("F9").Value = 402
$ws.Range("F10").Value = 3066
$ws.Range("F11").Value = 910
$ws.Range("G11").Value = 60
$ws.Range("F12").Value = 5232
$ws.Range("F13").Value = 1282
$ws.Range("F14").Value = 1031
$ws.Range("F16").Value = 845
$ws.Range("F18").Value = 406
$ws.Range("F19").Value = 1287
$ws.Range("F20").Value = 100
$ws.Range("F22").Value = 180
$ws.Range("F23").Value = 364
$ws.Range("F24").Value = 27
$ws.Range("F25").Value = 38
$ws.Range("F29").Value = 5077
$ws.Range("F31").Value = 4985
$ws.Range("F32").Value = 8982
$ws.Range("F34").Value = 152
$ws.Range("F35").Value = 136
$ws.Range("F36").Value = 223
$ws.Range("F37").Value = 429
$ws.Range("F39").Value = 85
$ws.Range("F40").Value = 4212
$ws.Range("F41").Value = 246

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 1133

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 572
$ws.Range("F5").Value = 42

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1344
$ws.Range("F4").Value = 385
$ws.Range("F5").Value = 327
$ws.Range("F6").Value = 881
$ws.Range("G6").Value = 60
$ws.Range("F7").Value = 760
$ws.Range("F8").Value = 402
$ws.Range("F9").Value = 910
$ws.Range("G9").Value = 60
$ws.Range("F11").Value = 1282
$ws.Range("F12").Value = 42
$ws.Range("F14").Value = 1031
$ws.Range("F16").Value = 845
$ws.Range("F17").Value = 406
$ws.Range("F19").Value = 1287
$ws.Range("F21").Value = 100
$ws.Range("F22").Value = 180
$ws.Range("F24").Value = 364
$ws.Range("F25").Value = 27
$ws.Range("F26").Value = 38
$ws.Range("F28").Value = 5077
$ws.Range("F30").Value = 8982
$ws.Range("F33").Value = 152
$ws.Range("F34").Value = 136
$ws.Range("F35").Value = 223
$ws.Range("F36").Value = 429
$ws.Range("F40").Value = 85
$ws.Range("F41").Value = 4212
$ws.Range("F48").Value = 246
